# cryptos.xlsx update - Tue Jan 16 11:09:39 UTC 2024 GitHub Actions run
#
# Refreshes the Price (D) / Volume(1h) (E) columns for every coin row, and
# swaps the two rows whose rank order changed (row 35 <-> row 36, i.e.
# WEMIXToken/LidoDAOToken trade places) by writing each destination cell
# with its final value directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.866.09"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "2.534.45"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.01"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.89"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0820"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("E13").Value = "  -3.59%  "
$ws.Range("D14").Value = "2.923.28"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "2.532.85"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").Value = "42.903.10"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.88"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.23%  "
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.69"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.33"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("E26").Value = "  -4.09%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.09"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.54"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.93"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.62"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.42"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.35"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.71"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("E39").Value = "  +10.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.119"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.83"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -12.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0305"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.82"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("D46").Value = "2.008.59"
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.24"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.51"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.28"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").Value = "2.777.90"
$ws.Range("E51").Value = "  -0.21%  "
